$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append at row 29 (matches columns A:K header layout)
$rowNum = 29

$values = @{
    "A" = "62,5"
    "B" = "110"
    "C" = "347"
    "D" = "406"
    "E" = "14682"
    "F" = "1"
    "G" = "86845"
    "H" = "2025-08-22 13:11"
    "I" = "A"
    "J" = "04"
    "K" = "CART.GRIS"
}

foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K")) {
    $cell = $ws.Range("$col$rowNum")
    # Force the cell to be stored as text so numeric-looking strings
    # (e.g. "110", "04") are not converted to numbers, matching the
    # inlineStr/text formatting used by the rest of the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col]
}
